$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data (shifts rows 1-3 down to 2-4)
$ws.Rows.Item(1).Insert()

# New header row (bold)
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Description(Optional)"
$ws.Range("D1").Value = "Placeholder"
$ws.Range("A1:D1").Font.Bold = $true

# Row 2 (discover/Discover) gains a description and placeholder
$ws.Range("C2").Value = "About Discover"
$ws.Range("D2").Value = "name:String "

# Row 3 (profile/Profile) gains a description
$ws.Range("C3").Value = "About Profile"

# Row 4 (formerly setting/Setting) is no longer part of the table
$ws.Rows.Item(4).Delete()

# Column widths for the new Description/Placeholder columns
$ws.Columns.Item(3).ColumnWidth = 17.6
$ws.Columns.Item(4).ColumnWidth = 16.5

# Match the saved selection from the authored workbook
$ws.Range("E6").Select()
